# Applies the "Updated who has which items" change to
# CIROH_Items_Comprehensive_List.xlsx
#
# Summary of the change:
#  - "2023-2024" sheet: view scrolled down (topLeftCell -> A23)
#  - "2024-2025" sheet: four new equipment rows (6-9) recording who
#    currently has the Microphone Setup, Sonar Sensor, PAR Sensor and
#    Prototype Controller, plus the new column G width and the
#    resulting selection/dimension bookkeeping.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "2023-2024" (first tab): the saved view had scrolled so row 23 is at
# the top of the window. (Best effort - some hosts don't persist window
# scroll position, but we still set it so the intent is recorded.)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("2023-2024")
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------------
# Sheet "2024-2025" (second tab, the active tab): add the four new rows.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2024-2025")
$ws2.Activate()

# The new block uses the same centered-text style as the existing rows.
$newBlock = $ws2.Range("A6:H9")
$newBlock.HorizontalAlignment = -4108   # xlCenter (center-aligned, matches existing data rows)

# Row 6 - Microphone Setup
$ws2.Range("A6").Value = "Microphone Setup"
$ws2.Range("B6").Value = 2
$ws2.Range("G6").Value = "Soheyl + Casey"
$ws2.Range("H6").Value = "Soheyl (x1)"

# Row 7 - Sonar Sensor
$ws2.Range("A7").Value = "Sonar Sensor"
$ws2.Range("B7").Value = 1
$ws2.Range("C7").Value = 99.99
$ws2.Range("C7").NumberFormat = "`"$`"#,##0.00_);[Red]\(`"$`"#,##0.00\)"
$ws2.Range("G7").Value = "Soheyl"
$ws2.Range("H7").Value = "Soheyl (x1)"

# Row 8 - PAR Sensor
$ws2.Range("A8").Value = "PAR Sensor"
$ws2.Range("B8").Value = 1
$ws2.Range("G8").Value = "Soheyl"
$ws2.Range("H8").Value = "Soheyl (x1)"

# Row 9 - Prototype Controller
$ws2.Range("A9").Value = "Prototype Controller"
$ws2.Range("B9").Value = 1
$ws2.Range("E9").Value = "DIY"
$ws2.Range("G9").Value = "Soheyl"
$ws2.Range("H9").Value = "Soheyl  (x1)"

# New column width introduced for column G alongside the new data.
$ws2.Columns.Item(7).ColumnWidth = 12.666666666666668

# Final saved selection on this sheet.
$ws2.Range("E13").Select()
